# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# de-de handback has completed (in addition to the already-recorded zh-cn
# handback), by:
#   - changing the "Status" text from "Ready for handoff" to
#     "Handed back: in sync with en-US" (affects Overview!E/F and the
#     Status column on both language sheets)
#   - filling in the "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns on the zh-cn and de-de sheets
#   - adding hyperlinks on the new "Latest Target File" cells, matching the
#     existing hyperlink on "Source File Name"

$wb = $excel.ActiveWorkbook

$urlC60e = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83aed6d0d3c6207ad315003c6446dd872a253475/e2e/c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"
$urlC8c5 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83aed6d0d3c6207ad315003c6446dd872a253475/e2e/c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: Status text is shared with the other sheets, updating
# the cells here updates the shared string used everywhere it is shown.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Row 2 (c60e3de8...)
$wsZh.Range("I2").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"
$wsZh.Range("I2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlC60e, "", "", "c60e3de8-f0b1-463f-83d6-957c38bb26a9.md") | Out-Null
$wsZh.Range("J2").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.6f9abbd882ba195e39018fd4a0cc57668314096a.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-02 02:55:58"

# Row 3 (c8c535ca...)
$wsZh.Range("I3").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"
$wsZh.Range("I3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlC8c5, "", "", "c8c535ca-e5aa-49e9-b056-1410c2db1b01.md") | Out-Null
$wsZh.Range("J3").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.206cdaf7b6d6b8c70b547477a8ef777d7d03ebc5.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-02 02:55:58"

# Widen the Status / Latest Target File / Latest Handback File columns to
# fit the new, longer content.
$wsZh.Range("C1").ColumnWidth = 29.1667
$wsZh.Range("I1").ColumnWidth = 39.1667
$wsZh.Range("J1").ColumnWidth = 39.1667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Row 2 (c60e3de8...)
$wsDe.Range("I2").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"
$wsDe.Range("I2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlC60e, "", "", "c60e3de8-f0b1-463f-83d6-957c38bb26a9.md") | Out-Null
$wsDe.Range("J2").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.6f9abbd882ba195e39018fd4a0cc57668314096a.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-02 02:56:11"

# Row 3 (c8c535ca...)
$wsDe.Range("I3").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"
$wsDe.Range("I3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlC8c5, "", "", "c8c535ca-e5aa-49e9-b056-1410c2db1b01.md") | Out-Null
$wsDe.Range("J3").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.206cdaf7b6d6b8c70b547477a8ef777d7d03ebc5.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-02 02:56:11"

# Widen the Status / Latest Target File / Latest Handback File columns to
# fit the new, longer content.
$wsDe.Range("C1").ColumnWidth = 29.1667
$wsDe.Range("I1").ColumnWidth = 39.1667
$wsDe.Range("J1").ColumnWidth = 39.1667

# ---------------------------------------------------------------------
# Overview sheet: widen the zh-cn / de-de status columns too, to match the
# new, longer Status text.
# ---------------------------------------------------------------------
$wsOverview.Range("E1").ColumnWidth = 29.1667
$wsOverview.Range("F1").ColumnWidth = 29.1667
